$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- n split (row 4): total stays 940, group split changes ---
$ws.Range("D4").Value = 644
$ws.Range("E4").Value = 296

# --- P-Value cells that look numeric must be forced to Text so they
#     remain stored as shared strings (matching the source table format) ---
foreach ($ref in @("F5","F6","F8","F10","F16","F26","F29")) {
    $ws.Range($ref).NumberFormat = "@"
}

# --- Updated summary statistics (D/E columns) and P-values (F column) ---
$ws.Range("D5").Value = "9.4 (6.4)"
$ws.Range("E5").Value = "9.3 (6.3)"
$ws.Range("F5").Value = "0.824"
$ws.Range("D6").Value = "311 (48.3)"
$ws.Range("E6").Value = "147 (49.7)"
$ws.Range("F6").Value = "0.749"
$ws.Range("D7").Value = "333 (51.7)"
$ws.Range("E7").Value = "149 (50.3)"
$ws.Range("D8").Value = "325 (50.5)"
$ws.Range("E8").Value = "141 (47.6)"
$ws.Range("F8").Value = "0.462"
$ws.Range("D9").Value = "319 (49.5)"
$ws.Range("E9").Value = "155 (52.4)"
$ws.Range("D10").Value = "477 (79.4)"
$ws.Range("E10").Value = "215 (78.2)"
$ws.Range("F10").Value = "0.141"
$ws.Range("D11").Value = "61 (10.1)"
$ws.Range("E11").Value = "41 (14.9)"
$ws.Range("D12").Value = "31 (5.2)"
$ws.Range("E12").Value = "12 (4.4)"
$ws.Range("D13").Value = "4 (0.7)"
$ws.Range("E13").Value = "1 (0.4)"
$ws.Range("D14").Value = "6 (1.0)"
$ws.Range("D15").Value = "22 (3.7)"
$ws.Range("E15").Value = "6 (2.2)"
$ws.Range("D16").Value = "124 (19.8)"
$ws.Range("E16").Value = "60 (21.1)"
$ws.Range("F16").Value = "0.730"
$ws.Range("D17").Value = "502 (80.2)"
$ws.Range("E17").Value = "225 (78.9)"
$ws.Range("D18").Value = "150 (26.5)"
$ws.Range("E18").Value = "108 (42.5)"
$ws.Range("F18").Value = "<0.001"
$ws.Range("D19").Value = "417 (73.5)"
$ws.Range("E19").Value = "146 (57.5)"
$ws.Range("D20").Value = "279 (43.3)"
$ws.Range("E20").Value = "184 (62.2)"
$ws.Range("F20").Value = "<0.001"
$ws.Range("D21").Value = "365 (56.7)"
$ws.Range("E21").Value = "112 (37.8)"
$ws.Range("D22").Value = "59.4 (24.7)"
$ws.Range("E22").Value = "73.7 (21.0)"
$ws.Range("D23").Value = "73 (11.5)"
$ws.Range("E23").Value = "55 (18.8)"
$ws.Range("D24").Value = "253 (39.9)"
$ws.Range("E24").Value = "200 (68.5)"
$ws.Range("D25").Value = "308 (48.6)"
$ws.Range("E25").Value = "37 (12.7)"
$ws.Range("D26").Value = "22 (3.4)"
$ws.Range("E26").Value = "14 (4.7)"
$ws.Range("F26").Value = "0.487"
$ws.Range("D27").Value = "344 (53.4)"
$ws.Range("E27").Value = "163 (55.1)"
$ws.Range("D28").Value = "278 (43.2)"
$ws.Range("E28").Value = "119 (40.2)"
$ws.Range("D29").Value = "108 (16.8)"
$ws.Range("E29").Value = "56 (18.9)"
$ws.Range("F29").Value = "0.488"
$ws.Range("D30").Value = "534 (83.2)"
$ws.Range("E30").Value = "240 (81.1)"
